# Burndown - Sprint 2: log today's time spent / remaining work.
# (Meeting Recording & Burndown Chart)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# "Room Allocation System" row - estimate bumped from 2 to 3 hours.
$ws.Range("B12").Value = 3

# Day 3 (column E) hours logged against several tasks.
$ws.Range("E12").Value = 2      # Room Allocation System
$ws.Range("E13").Value = 1      # AI Patrolling
$ws.Range("E14").Value = 3      # Basement Room A: Event
$ws.Range("E19").Value = 4      # Basement Room F: Event
$ws.Range("E21").Value = 0.05   # User Stories (Final Event)

$ws.Range("E12").Select()
